$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Poisson"
$ws.Range("C2").Value = "FE"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "Regular_RightTriangles"
$ws.Range("F2").Value = 2.003941254089581
$ws.Range("G2").Value = "Triangles"
$ws.Range("H2").Value = "Green"
$ws.Range("I2").Value = 119.5647480487823

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Poisson"
$ws.Range("C3").Value = "FE"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "Unstructured_triangles"
$ws.Range("F3").Value = 2.015608644460145
$ws.Range("G3").Value = "Triangles"
$ws.Range("H3").Value = "Green"
$ws.Range("I3").Value = 6.094477891921997

$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "Poisson"
$ws.Range("C4").Value = "FE"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "Regular_Tetrahedra"
$ws.Range("F4").Value = 1.340336836145038
$ws.Range("G4").Value = "Tetrahedron"
$ws.Range("H4").Value = "Green"
$ws.Range("I4").Value = 208.4592311382294

$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Poisson"
$ws.Range("C5").Value = "FE"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = "Unstructured_Tetrahedra"
$ws.Range("F5").Value = 0.6690820358074518
$ws.Range("G5").Value = "Tetrahedron"
$ws.Range("H5").Value = "Green"
$ws.Range("I5").Value = 11.90931582450867

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Poisson"
$ws.Range("C6").Value = "FV"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "RegularSquares"
$ws.Range("F6").Value = 2.003941213535303
$ws.Range("G6").Value = "Squares"
$ws.Range("H6").Value = "Green"
$ws.Range("I6").Value = 9.851321935653687

$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "Poisson"
$ws.Range("C7").Value = "FV"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "RegularSquares"
$ws.Range("F7").Value = 2.003941211551183
$ws.Range("G7").Value = "Squares"
$ws.Range("H7").Value = "Green"
$ws.Range("I7").Value = 9.898462057113647

$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Poisson"
$ws.Range("C8").Value = "FV"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "Regular_RightTriangles"
$ws.Range("F8").Value = 0.02119663100406134
$ws.Range("G8").Value = "Triangles"
$ws.Range("H8").Value = "Green"
$ws.Range("I8").Value = 15.68182492256165

$ws.Range("A9").Value = 0
$ws.Range("B9").Value = "Poisson"
$ws.Range("C9").Value = "FV"
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = "Regular_RightTriangles"
$ws.Range("F9").Value = -0.005617740418916485
$ws.Range("G9").Value = "Triangles"
$ws.Range("H9").Value = "Orange(order 0)"
$ws.Range("I9").Value = 15.86531114578247

$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Poisson"
$ws.Range("C10").Value = "FV"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = "Structured_triangles"
$ws.Range("F10").Value = 0.8952237869134417
$ws.Range("G10").Value = "Triangles"
$ws.Range("H10").Value = "Green"
$ws.Range("I10").Value = 4.744688034057617

$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Poisson"
$ws.Range("C11").Value = "FV"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = "Unstructured_triangles"
$ws.Range("F11").Value = 0.6137798580984465
$ws.Range("G11").Value = "Triangles"
$ws.Range("H11").Value = "Green"
$ws.Range("I11").Value = 2.505467891693115

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Poisson"
$ws.Range("C12").Value = "FV"
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = "Regular_Cubes"
$ws.Range("F12").Value = 1.340336836132099
$ws.Range("G12").Value = "Cubes"
$ws.Range("H12").Value = "Green"
$ws.Range("I12").Value = 5.868787050247192

$ws.Range("A13").Value = 9
$ws.Range("B13").Value = "Poisson"
$ws.Range("C13").Value = "FV"
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = "Regular_Tetrahedra"
$ws.Range("F13").Value = 0.006535470643459771
$ws.Range("G13").Value = "Tetrahedron"
$ws.Range("H13").Value = "Green"
$ws.Range("I13").Value = 62.12639307975769

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Poisson"
$ws.Range("C14").Value = "FV"
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = "Unstructured_Tetrahedra"
$ws.Range("F14").Value = 0.5358788100873692
$ws.Range("G14").Value = "Tetrahedron"
$ws.Range("H14").Value = "Green"
$ws.Range("I14").Value = 3.679362773895264
